$d = $word.ActiveDocument

# Helper: find the first occurrence of $find in the document and replace
# it in-place by assigning straight to the found Range's .Text property.
# (Passing the replacement through Find.Execute's own Replace argument
# triggers this runtime's "smart quotes" autocorrect, which mangles the
# straight apostrophes used throughout this document - so we avoid that
# path and set Range.Text directly on the located match instead.)
function ReplaceText($find, $replace) {
    $rng = $d.Content
    $found = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("Find failed for: " + $find)
    }
    $rng.Text = $replace
}

# Helper: set a table cell's text. Tables/Cells are re-resolved fresh
# each call (instead of reusing a previously fetched Cell/Table object)
# because earlier edits shift character offsets and stale cell handles
# end up pointing at the wrong cell.
function SetCellText($tableIndex, $row, $col, $text) {
    $d.Tables.Item($tableIndex).Cell($row, $col).Range.Text = $text
}

# --- Table 1 (resourcesResponse) ---

# caseId description: "Identifiant..." -> "A valoriser avec l'identifiant..."
ReplaceText "Identifiant partagé de l'affaire/dossier, généré une seule fois" "A valoriser avec l'identifiant partagé de l'affaire/dossier, généré une seule fois"

# requestId description: "Identifiant..." -> "A valoriser avec l'identifiant..."
ReplaceText "Identifiant unique partagé de la demande de ressource,  généré une seule fois" "A valoriser avec l'identifiant unique partagé de la demande de ressource,  généré une seule fois"

# --- Swap content between "response" row (table 1, row 4) and
#     "datetime" row (table 2, row 2): the old "response" Description/
#     Exemple text moves down to "datetime", and "response" gets new
#     Description text with an emptied Exemple cell. ---

# Move the old "response" row text down onto the "datetime" row first.
SetCellText 2 2 5 "Groupe date heure de début de la demande"
SetCellText 2 2 6 "2022-09-27T08:23:34+02:00"

# Give the "response" row its new description and clear its example.
SetCellText 1 4 5 "Objet permettant de transmettre les détails de la réponse à une demande de ressource"
SetCellText 1 4 6 ""

# --- Table 2 (response) remaining edits ---

# answer row: Description (two lines, joined with a line break) + Exemple
$answerDesc = "A valoriser avec la réponse apportée. Cf Nomenclature associée" + [char]11 + "ACCEPTEE, REFUSEE, PARTIELLE, DIFFEREE"
SetCellText 2 3 5 $answerDesc
SetCellText 2 3 6 "ACCEPTEE"

# deadline row: Description + Exemple
ReplaceText "Indique le délai de réponse auquel s'engage l'expéditeur (en minutes)" "A valoriser avec le délai de réponse auquel s'engage l'expéditeur (en minutes),"
SetCellText 2 4 6 "10"

# freetext row: Description + Exemple
ReplaceText "Commentaire libre pour apporter toutes précisions utiles à la réponse." "Commentaire libre permettant d'apporter toutes précisions utiles à la réponse."
SetCellText 2 5 6 "SMUR 1 non dispo"
